$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(164).Insert()

$ws.Range("A164").Value = 3
$ws.Range("B164").Value = "Femacal de La Calera"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44825
$ws.Range("E164").Value = 5
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100101
$ws.Range("H164").Value = "Berries"
$ws.Range("I164").Value = 100101001
$ws.Range("J164").Value = "Arándano (blue)"
$ws.Range("K164").Value = "Sin especificar"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 56
$ws.Range("N164").Value = 12000
$ws.Range("O164").Value = 12000
$ws.Range("P164").Value = 12000
$ws.Range("Q164").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R164").Value = "Provincia de Limarí"
$ws.Range("S164").Value = 8000
$ws.Range("T164").Value = 1.5
